$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D21").Value = -6932.12
$ws.Range("F21").Value = -1540.471111111111
$ws.Range("H21").Value = -51270.12
$ws.Range("I21").Value = -11233.47111111111

$ws.Range("D22").Value = 25000
$ws.Range("F22").Value = 5555.555555555556
$ws.Range("H22").Value = 1116
$ws.Range("I22").Value = 1445.555555555556

$ws.Range("D23").Value = 25000
$ws.Range("F23").Value = 4545.454545454545
$ws.Range("H23").Value = 2326
$ws.Range("I23").Value = 4393.454545454545

$ws.Range("D24").Value = 25000
$ws.Range("F24").Value = 4545.454545454545
$ws.Range("H24").Value = 12930
$ws.Range("I24").Value = 2884.454545454545

$ws.Range("D25").Value = -6932.12
$ws.Range("F25").Value = -1733.03
$ws.Range("H25").Value = -27076.12
$ws.Range("I25").Value = -4351.03

$ws.Range("D26").Value = -6932.12
$ws.Range("F26").Value = -1540.471111111111
$ws.Range("H26").Value = -37025.12
$ws.Range("I26").Value = -8251.47111111111

$ws.Range("D35").Value = -6932.12
$ws.Range("F35").Value = -1260.385454545454
$ws.Range("H35").Value = -26686.12
$ws.Range("I35").Value = -5294.385454545454

$ws.Range("D36").Value = -6932.12
$ws.Range("F36").Value = -1733.03
$ws.Range("H36").Value = -11793.12
$ws.Range("I36").Value = -3182.03

$ws.Range("D37").Value = -6932.12
$ws.Range("F37").Value = -1540.471111111111
$ws.Range("H37").Value = -21756.12
$ws.Range("I37").Value = -5093.471111111111

$ws.Range("D350").Value = -6932.12
$ws.Range("F350").Value = -1540.471111111111
$ws.Range("H350").Value = -20729.12
$ws.Range("I350").Value = -3765.471111111111

$ws.Range("D351").Value = -6932.12
$ws.Range("F351").Value = -1540.471111111111
$ws.Range("H351").Value = -13263.12
$ws.Range("I351").Value = -2706.471111111111

$ws.Range("D352").Value = -6932.12
$ws.Range("F352").Value = -1540.471111111111
$ws.Range("H352").Value = -28786.12
$ws.Range("I352").Value = -5002.471111111111

$ws.Range("D354").Value = -6932.12
$ws.Range("F354").Value = -1540.471111111111
$ws.Range("H354").Value = -9097.119999999999
$ws.Range("I354").Value = -2187.471111111111

$ws.Range("D355").Value = -6932.12
$ws.Range("F355").Value = -1540.471111111111
$ws.Range("H355").Value = -10805.12
$ws.Range("I355").Value = -2814.471111111111

$ws.Range("D356").Value = -6932.12
$ws.Range("F356").Value = -1733.03
$ws.Range("H356").Value = -26509.12
$ws.Range("I356").Value = -4755.03

$ws.Range("D357").Value = -6932.12
$ws.Range("F357").Value = -1540.471111111111
$ws.Range("H357").Value = -19425.12
$ws.Range("I357").Value = -2533.471111111111

$ws.Range("D358").Value = -6932.12
$ws.Range("F358").Value = -1733.03
$ws.Range("H358").Value = -28212.12
$ws.Range("I358").Value = -6404.03

$ws.Range("D359").Value = -6932.12
$ws.Range("F359").Value = -1260.385454545454
$ws.Range("H359").Value = -33075.12
$ws.Range("I359").Value = -7359.385454545454
